$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("O7").Value = 1.14
$ws.Range("P7").Value = 5.86
$ws.Range("S7").Value = 1.26
$ws.Range("T7").Value = 3.92

# Row 12
$ws.Range("G12").Value = 1.2
$ws.Range("H12").Value = 5.75
$ws.Range("J12").Value = 1.67
$ws.Range("K12").Value = 2.5
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 13
$ws.Range("O12").Value = 1.25
$ws.Range("P12").Value = 3.75
$ws.Range("Q12").Value = 1.85
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 1.36
$ws.Range("T12").Value = 3
$ws.Range("W12").Value = 5.5
$ws.Range("AC12").Value = 9.5
$ws.Range("AD12").Value = 12
$ws.Range("AH12").Value = 29
$ws.Range("AT12").Value = 3

# Row 23
$ws.Range("G23").Value = 2.62
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 2.67
$ws.Range("J23").Value = 3.25
$ws.Range("K23").Value = 2
$ws.Range("L23").Value = 3.2
$ws.Range("M23").Value = 1.09
$ws.Range("N23").Value = 7
$ws.Range("O23").Value = 1.37
$ws.Range("P23").Value = 2.62
$ws.Range("Q23").Value = 2.07
$ws.Range("R23").Value = 1.6
$ws.Range("S23").Value = 1.42
$ws.Range("T23").Value = 2.47
$ws.Range("U23").Value = 1.8
$ws.Range("V23").Value = 1.8
$ws.Range("W23").Value = 7.3
$ws.Range("Y23").Value = 10
$ws.Range("AA23").Value = 25
$ws.Range("AB23").Value = 37
$ws.Range("AC23").Value = 7.7
$ws.Range("AD23").Value = 5.8
$ws.Range("AE23").Value = 14.5
$ws.Range("AF23").Value = 75
$ws.Range("AG23").Value = 700
$ws.Range("AH23").Value = 7.8
$ws.Range("AI23").Value = 13
$ws.Range("AO23").Value = 14.5
$ws.Range("AP23").Value = 22
$ws.Range("AQ23").Value = 65
$ws.Range("AR23").Value = 100
$ws.Range("AS23").Value = 300
$ws.Range("AT23").Value = 2.42
$ws.Range("AU23").Value = 6.8
$ws.Range("AV23").Value = 60
$ws.Range("AW23").Value = 4.5
$ws.Range("AX23").Value = 14

# Row 38
$ws.Range("H38").Value = 3.25
$ws.Range("J38").Value = 4.25
$ws.Range("K38").Value = 2.05
$ws.Range("L38").Value = 2.6
$ws.Range("W38").Value = 9.5
$ws.Range("Y38").Value = 13
$ws.Range("AB38").Value = 50
$ws.Range("AH38").Value = 6.1
$ws.Range("AI38").Value = 8.25
$ws.Range("AL38").Value = 17.5
$ws.Range("AM38").Value = 35
$ws.Range("AO38").Value = 22
$ws.Range("AR38").Value = 175
$ws.Range("AU38").Value = 7.8
$ws.Range("AV38").Value = 80
$ws.Range("AX38").Value = 10.25
$ws.Range("AY38").Value = 22
$ws.Range("BA38").Value = 90

